$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("caso abc imptos 10%")

# Rename the first sheet
$ws.Name = "BaseCase"

# Update cell values (accented -> non-accented versions)
$ws.Range("A4").Value = "Inversion"
$ws.Range("A10").Value = "Costes operacion"
$ws.Range("A13").Value = "Amortizacion"

# Update the view: scroll and selection
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B4").Select()
